$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = 322348.0527887827
$ws.Range("D2").Value = 0.2243460764587525
$ws.Range("E2").Value = 0.2234468937875752
$ws.Range("F2").Value = 0.2238955823293172
$ws.Range("G2").Value = 0.0006904096578525566
$ws.Range("C3").Value = 51611499.46010433
$ws.Range("D3").Value = 0.3856800638807559
$ws.Range("E3").Value = 0.386296987470008
$ws.Range("F3").Value = 0.3859882791688865
$ws.Range("G3").Value = 0.00005619514953090941
$ws.Range("C4").Value = 25855320.12366612
$ws.Range("D4").Value = 0.393186372745491
$ws.Range("E4").Value = 0.3922953878965609
$ws.Range("F4").Value = 0.3927403749916595
$ws.Range("G4").Value = 0.0001136965890482947
$ws.Range("C5").Value = 51650891.15947316
$ws.Range("D5").Value = 0.7949953164726348
$ws.Range("E5").Value = 0.7919221540922421
$ws.Range("F5").Value = 0.7934557595993321
$ws.Range("G5").Value = 0.0001147994692517187
$ws.Range("C6").Value = 40844390.64762544
$ws.Range("D6").Value = 0.3856800638807559
$ws.Range("E6").Value = 0.386296987470008
$ws.Range("F6").Value = 0.3859882791688865
$ws.Range("G6").Value = 0.00007100891661468888
$ws.Range("C7").Value = 25874329.39525799
$ws.Range("E7").Value = 0.7980538523060517
$ws.Range("F7").Value = 0.8876862628808658
$ws.Range("G7").Value = 0.0002053996289017543
$ws.Range("C8").Value = 40875555.94910018
$ws.Range("D8").Value = 0.7949953164726348
$ws.Range("E8").Value = 0.7919221540922421
$ws.Range("F8").Value = 0.7934557595993321
$ws.Range("G8").Value = 0.0001450621221854314
$ws.Range("C9").Value = 20467034.76685503
$ws.Range("D9").Value = 0.393186372745491
$ws.Range("E9").Value = 0.3922953878965609
$ws.Range("F9").Value = 0.3927403749916595
$ws.Range("G9").Value = 0.000143629096266214
$ws.Range("C10").Value = 20482074.1853513
$ws.Range("E10").Value = 0.7980538523060517
$ws.Range("F10").Value = 0.8876862628808658
$ws.Range("G10").Value = 0.0002594745828851996
